$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price / volume(1h) data (and a few reordered rows)
# to match the latest scrape, per the automated GitHub Actions update job.

$ws.Range("D2").Value = "66.208.64"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").Value = "3.519.30"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'583.34"
$ws.Range("E5").Value = "  +6.01%  "

$ws.Range("D6").Value = "'179.54"
$ws.Range("E6").Value = "  -4.88%  "

$ws.Range("D7").Value = "'0.634"
$ws.Range("E7").Value = "  +4.80%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +1.58%  "

$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  +6.48%  "

$ws.Range("D11").Value = "'55.91"
$ws.Range("E11").Value = "  +2.40%  "

$ws.Range("D12").Value = "'0.0000280"
$ws.Range("E12").Value = "  +4.22%  "

$ws.Range("D13").Value = "'9.31"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("D14").Value = "4.084.09"
$ws.Range("E14").Value = "  -0.06%  "

$ws.Range("D15").Value = "3.525.83"
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("E16").Value = "  +0.63%  "

$ws.Range("D17").Value = "'18.43"
$ws.Range("E17").Value = "  +1.60%  "

$ws.Range("D18").Value = "'12.12"
$ws.Range("E18").Value = "  +2.26%  "

$ws.Range("D19").Value = "66.234.09"
$ws.Range("E19").Value = "  -1.26%  "

$ws.Range("E20").Value = "  +2.25%  "

$ws.Range("D21").Value = "'415.55"
$ws.Range("E21").Value = "  -3.26%  "

$ws.Range("E22").Value = "  +10.92%  "

$ws.Range("D23").Value = "'4.46"
$ws.Range("E23").Value = "  +7.43%  "

$ws.Range("D24").Value = "'85.38"
$ws.Range("E24").Value = "  +0.53%  "

$ws.Range("D25").Value = "'13.59"
$ws.Range("E25").Value = "  +13.28%  "

$ws.Range("D26").Value = "'11.16"
$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("D27").Value = "'2.88"
$ws.Range("E27").Value = "  -0.66%  "

$ws.Range("D28").Value = "'6.06"
$ws.Range("E28").Value = "  -1.48%  "

$ws.Range("D29").Value = "'9.22"
$ws.Range("E29").Value = "  +2.88%  "

$ws.Range("D30").Value = "'30.43"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").Value = "'6.65"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").Value = "'610.13"
$ws.Range("E32").Value = "  -5.44%  "

$ws.Range("D33").Value = "'11.82"
$ws.Range("E33").Value = "  +1.12%  "

$ws.Range("E34").Value = "  +0.88%  "

$ws.Range("D35").Value = "'60.25"
$ws.Range("E35").Value = "  +1.61%  "

$ws.Range("D36").Value = "'0.155"
$ws.Range("E36").Value = "  +9.41%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'3.67"
$ws.Range("E38").Value = "  +9.63%  "

$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0801"
$ws.Range("E39").Value = "  -1.10%  "

$ws.Range("D40").Value = "'36.96"
$ws.Range("E40").Value = "  -3.90%  "

$ws.Range("D41").Value = "'0.387"
$ws.Range("E41").Value = "  -0.65%  "

$ws.Range("D42").Value = "3.259.50"
$ws.Range("E42").Value = "  +7.84%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "'3.01"
$ws.Range("E44").Value = "  +5.30%  "

$ws.Range("D45").Value = "'3.40"
$ws.Range("E45").Value = "  +3.10%  "

$ws.Range("E46").Value = "  -2.72%  "

$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("E48").Value = "  -7.25%  "

$ws.Range("E49").Value = "  +2.25%  "

$ws.Range("D50").Value = "'8.68"
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("D51").Value = "'139.91"
$ws.Range("E51").Value = "  -0.76%  "

